$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fecha" column (G) held values like "2002-1" .. "2018-4". Re-key them
# to the "2002T1" .. "2018T4" style (hyphen -> "T") - this is what moved the
# underlying shared-string table around in the saved file.
for ($r = 2; $r -le 69; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $old = $cell.Value2
    $new = $old -replace "-", "T"
    $cell.Value = $new
}

# Column A ("INVERSION") picks up an explicit custom width (no longer
# auto "best fit"), and the newly-used column E ("Trimestre") gets its own
# custom width too.
$ws.Columns.Item(1).ColumnWidth = 10.3
$ws.Columns.Item(5).ColumnWidth = 9

# Scroll the view back to the top and move the selection to F3 (previously
# it was parked at B70 with the view scrolled down to row 67).
$ws.Range("A1").Select()
$ws.Range("F3").Select()
